# Updated cryptos list on Sun Aug 11 11:20:49 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "60.673.24"
$ws.Cells.Item(2,5).Value = "  -0.13%  "
$ws.Cells.Item(3,4).Value = "2.649.09"
$ws.Cells.Item(3,5).Value = "  +1.16%  "
$ws.Cells.Item(4,5).Value = "  +0.04%  "
$cell = $ws.Cells.Item(5,4)
$cell.Value = "'525.95"
$cell.Style = "Normal"
$ws.Cells.Item(5,5).Value = "  +2.10%  "
$cell = $ws.Cells.Item(6,4)
$cell.Value = "'154.12"
$cell.Style = "Normal"
$ws.Cells.Item(6,5).Value = "  -0.29%  "
$cell = $ws.Cells.Item(7,4)
$cell.Value = "'0.998"
$cell.Style = "Normal"
$ws.Cells.Item(7,5).Value = "  +0.12%  "
$cell = $ws.Cells.Item(8,4)
$cell.Value = "'0.580"
$cell.Style = "Normal"
$ws.Cells.Item(8,5).Value = "  -1.52%  "
$ws.Cells.Item(9,5).Value = "  -3.81%  "
$ws.Cells.Item(10,5).Value = "  +3.89%  "
$ws.Cells.Item(11,5).Value = "  +0.49%  "
$ws.Cells.Item(12,5).Value = "  -0.68%  "
$ws.Cells.Item(13,4).Value = "3.112.29"
$ws.Cells.Item(13,5).Value = "  +1.20%  "
$ws.Cells.Item(14,4).Value = "60.707.65"
$ws.Cells.Item(14,5).Value = "  -0.03%  "
$cell = $ws.Cells.Item(15,4)
$cell.Value = "'21.89"
$cell.Style = "Normal"
$ws.Cells.Item(15,5).Value = "  +0.52%  "
$ws.Cells.Item(16,5).Value = "  +0.37%  "
$ws.Cells.Item(17,4).Value = "2.648.19"
$ws.Cells.Item(17,5).Value = "  +0.89%  "
$cell = $ws.Cells.Item(18,4)
$cell.Value = "'4.75"
$cell.Style = "Normal"
$ws.Cells.Item(18,5).Value = "  -0.24%  "
$cell = $ws.Cells.Item(19,4)
$cell.Value = "'353.06"
$cell.Style = "Normal"
$ws.Cells.Item(19,5).Value = "  -1.33%  "
$cell = $ws.Cells.Item(20,4)
$cell.Value = "'10.62"
$cell.Style = "Normal"
$ws.Cells.Item(20,5).Value = "  -0.64%  "
$cell = $ws.Cells.Item(21,4)
$cell.Value = "'6.28"
$cell.Style = "Normal"
$ws.Cells.Item(21,5).Value = "  +1.22%  "
$cell = $ws.Cells.Item(22,4)
$cell.Value = "'1.00"
$cell.Style = "Normal"
$ws.Cells.Item(22,5).Value = "  +0.34%  "
$cell = $ws.Cells.Item(23,4)
$cell.Value = "'61.40"
$cell.Style = "Normal"
$ws.Cells.Item(23,5).Value = "  +0.89%  "
$cell = $ws.Cells.Item(24,4)
$cell.Value = "'0.428"
$cell.Style = "Normal"
$ws.Cells.Item(24,5).Value = "  +0.33%  "
$ws.Cells.Item(25,5).Value = "  -0.02%  "
$ws.Cells.Item(26,5).Value = "  +0.36%  "
$ws.Cells.Item(27,4).Value = "0.0₃0848"
$ws.Cells.Item(27,5).Value = "  -0.03%  "
$cell = $ws.Cells.Item(28,4)
$cell.Value = "'7.26"
$cell.Style = "Normal"
$ws.Cells.Item(28,5).Value = "  -1.78%  "
$ws.Cells.Item(29,5).Value = "  +0.01%  "
$cell = $ws.Cells.Item(30,4)
$cell.Value = "'6.16"
$cell.Style = "Normal"
$ws.Cells.Item(30,5).Value = "  +2.94%  "
$cell = $ws.Cells.Item(31,4)
$cell.Value = "'19.39"
$cell.Style = "Normal"
$ws.Cells.Item(31,5).Value = "  -0.63%  "
$cell = $ws.Cells.Item(32,4)
$cell.Value = "'1.61"
$cell.Style = "Normal"
$ws.Cells.Item(32,5).Value = "  +1.75%  "
$cell = $ws.Cells.Item(33,4)
$cell.Value = "'149.95"
$cell.Style = "Normal"
$ws.Cells.Item(33,5).Value = "  -1.18%  "
$cell = $ws.Cells.Item(34,4)
$cell.Value = "'4.09"
$cell.Style = "Normal"
$ws.Cells.Item(34,5).Value = "  +1.32%  "
$ws.Cells.Item(35,5).Value = "  -0.90%  "
$cell = $ws.Cells.Item(36,4)
$cell.Value = "'0.909"
$cell.Style = "Normal"
$ws.Cells.Item(36,5).Value = "  +6.13%  "
$cell = $ws.Cells.Item(37,4)
$cell.Value = "'0.885"
$cell.Style = "Normal"
$ws.Cells.Item(37,5).Value = "  -0.44%  "
$cell = $ws.Cells.Item(38,4)
$cell.Value = "'36.83"
$cell.Style = "Normal"
$ws.Cells.Item(38,5).Value = "  +1.13%  "
$cell = $ws.Cells.Item(39,4)
$cell.Value = "'305.42"
$cell.Style = "Normal"
$ws.Cells.Item(39,5).Value = "  +3.82%  "
$ws.Cells.Item(40,5).Value = "  -1.87%  "
$cell = $ws.Cells.Item(41,4)
$cell.Value = "'3.75"
$cell.Style = "Normal"
$ws.Cells.Item(41,5).Value = "  -0.28%  "
$cell = $ws.Cells.Item(42,4)
$cell.Value = "'0.639"
$cell.Style = "Normal"
$ws.Cells.Item(42,5).Value = "  +2.49%  "
$ws.Cells.Item(43,5).Value = "  +0.15%  "
$cell = $ws.Cells.Item(44,4)
$cell.Value = "'20.21"
$cell.Style = "Normal"
$ws.Cells.Item(44,5).Value = "  +1.87%  "
$cell = $ws.Cells.Item(45,4)
$cell.Value = "'0.0560"
$cell.Style = "Normal"
$ws.Cells.Item(45,5).Value = "  +0.39%  "
$cell = $ws.Cells.Item(46,4)
$cell.Value = "'0.998"
$cell.Style = "Normal"
$ws.Cells.Item(46,5).Value = "  +0.21%  "
$cell = $ws.Cells.Item(47,4)
$cell.Value = "'0.0240"
$cell.Style = "Normal"
$ws.Cells.Item(47,5).Value = "  +1.74%  "
$cell = $ws.Cells.Item(48,4)
$cell.Value = "'4.84"
$cell.Style = "Normal"
$ws.Cells.Item(48,5).Value = "  -2.58%  "
$ws.Cells.Item(49,2).Value = "WhiteBITCoin"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$cell = $ws.Cells.Item(49,4)
$cell.Value = "'10.36"
$cell.Style = "Normal"
$ws.Cells.Item(49,5).Value = "  +0.50%  "
$ws.Cells.Item(50,2).Value = "InjectiveProtocol"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell = $ws.Cells.Item(50,4)
$cell.Value = "'19.16"
$cell.Style = "Normal"
$ws.Cells.Item(50,5).Value = "  +2.77%  "
$ws.Cells.Item(51,4).Value = "1.982.10"
$ws.Cells.Item(51,5).Value = "  -0.78%  "

Write-Output "Updated cryptos list with latest price and volume data"
